# Auto-generated update of cached market-price / profit figures in the
# FFXIV "Shinryu_Profits" workbook (per-job profit tables), as produced by
# the scheduled data-refresh runner. Each worksheet (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) holds a Leve profit table; columns H-N are the
# recalculated market price / profit figures for specific rows.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 162.47368
$ws.Range("I33").Value = 38.785713
$ws.Range("J33").Value = 508.8
$ws.Range("K33").Value = 38.785713
$ws.Range("L33").Value = 508.8
$ws.Range("M33").Value = 190.214287
$ws.Range("N33").Value = -966.8
# Row 40
$ws.Range("H40").Value = 1762.8422
$ws.Range("I40").Value = 1775.125
$ws.Range("J40").Value = 1697.3334
$ws.Range("K40").Value = 1775.125
$ws.Range("L40").Value = 1697.3334
$ws.Range("M40").Value = -1600.125
$ws.Range("N40").Value = -2047.3334
# Row 43
$ws.Range("H43").Value = 784.3226
$ws.Range("I43").Value = 662.8570999999999
$ws.Range("J43").Value = 819.75
$ws.Range("K43").Value = 662.8570999999999
$ws.Range("L43").Value = 819.75
$ws.Range("M43").Value = -593.8570999999999
$ws.Range("N43").Value = -957.75
# Row 76
$ws.Range("H76").Value = 3587.6
$ws.Range("I76").Value = 3396
$ws.Range("K76").Value = 3396
$ws.Range("M76").Value = -3081
# Row 79
$ws.Range("H79").Value = 3587.6
$ws.Range("I79").Value = 3396
$ws.Range("K79").Value = 3396
$ws.Range("M79").Value = -2304
# Row 93
$ws.Range("H93").Value = 90114.36
$ws.Range("J93").Value = 90114.36
$ws.Range("L93").Value = 90114.36
$ws.Range("N93").Value = -95106.36
# Row 97
$ws.Range("H97").Value = 983.3333
$ws.Range("J97").Value = 983.3333
$ws.Range("L97").Value = 2949.9999
$ws.Range("N97").Value = -3941.9999
# Row 101
$ws.Range("H101").Value = 6099.4165
$ws.Range("I101").Value = 417.84616
$ws.Range("J101").Value = 12814
$ws.Range("K101").Value = 1253.53848
$ws.Range("L101").Value = 38442
$ws.Range("M101").Value = 368.4615200000001
$ws.Range("N101").Value = -41686
# Row 112
$ws.Range("H112").Value = 1428.8823
$ws.Range("I112").Value = 539.8
$ws.Range("J112").Value = 1799.3334
$ws.Range("K112").Value = 1619.4
$ws.Range("L112").Value = 5398.0002
$ws.Range("M112").Value = -511.3999999999999
$ws.Range("N112").Value = -7614.0002
# Row 132
$ws.Range("H132").Value = 3737.8708
$ws.Range("I132").Value = 3549.7
$ws.Range("J132").Value = 4080
$ws.Range("K132").Value = 10649.1
$ws.Range("L132").Value = 12240
$ws.Range("M132").Value = -8119.099999999999
$ws.Range("N132").Value = -17300
# Row 137
$ws.Range("H137").Value = 2476.875
$ws.Range("J137").Value = 4183.0835
$ws.Range("L137").Value = 12549.2505
$ws.Range("N137").Value = -17649.2505

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2217
$ws.Range("I61").Value = 1222.7
$ws.Range("K61").Value = 1222.7
$ws.Range("M61").Value = -1010.7
# Row 74
$ws.Range("H74").Value = 1322.075
$ws.Range("I74").Value = 1343.2759
$ws.Range("K74").Value = 1343.2759
$ws.Range("M74").Value = -469.2759000000001
# Row 77
$ws.Range("H77").Value = 1322.075
$ws.Range("I77").Value = 1343.2759
$ws.Range("K77").Value = 6716.379500000001
$ws.Range("M77").Value = -2348.379500000001
# Row 110
$ws.Range("H110").Value = 1678.7084
$ws.Range("I110").Value = 1381.8823
$ws.Range("J110").Value = 2399.5715
$ws.Range("K110").Value = 1381.8823
$ws.Range("L110").Value = 2399.5715
$ws.Range("M110").Value = 663.1177
$ws.Range("N110").Value = -6489.5715
# Row 122
$ws.Range("H122").Value = 2188.889
$ws.Range("I122").Value = 2099.8
$ws.Range("J122").Value = 2300.25
$ws.Range("K122").Value = 6299.400000000001
$ws.Range("L122").Value = 6900.75
$ws.Range("M122").Value = -3849.400000000001
$ws.Range("N122").Value = -11800.75
# Row 136
$ws.Range("H136").Value = 2217
$ws.Range("I136").Value = 1222.7
$ws.Range("K136").Value = 3668.1
$ws.Range("M136").Value = -1118.1

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 54
$ws.Range("H54").Value = 2332.6667
$ws.Range("I54").Value = 591.5
$ws.Range("J54").Value = 9297.333000000001
$ws.Range("K54").Value = 591.5
$ws.Range("L54").Value = 9297.333000000001
$ws.Range("M54").Value = -107.5
$ws.Range("N54").Value = -10265.333
# Row 99
$ws.Range("H99").Value = 1717.2858
$ws.Range("I99").Value = 1112.8572
$ws.Range("K99").Value = 1112.8572
$ws.Range("M99").Value = 385.1428000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1995.3077
$ws.Range("I31").Value = 1560.2727
$ws.Range("J31").Value = 4388
$ws.Range("K31").Value = 1560.2727
$ws.Range("L31").Value = 4388
$ws.Range("M31").Value = -1265.2727
$ws.Range("N31").Value = -4978
# Row 34
$ws.Range("H34").Value = 1995.3077
$ws.Range("I34").Value = 1560.2727
$ws.Range("J34").Value = 4388
$ws.Range("K34").Value = 1560.2727
$ws.Range("L34").Value = 4388
$ws.Range("M34").Value = -1358.2727
$ws.Range("N34").Value = -4792
# Row 58
$ws.Range("H58").Value = 1011.1316
$ws.Range("I58").Value = 796.5333000000001
$ws.Range("J58").Value = 1815.875
$ws.Range("K58").Value = 796.5333000000001
$ws.Range("L58").Value = 1815.875
$ws.Range("M58").Value = -593.5333000000001
$ws.Range("N58").Value = -2221.875
# Row 107
$ws.Range("H107").Value = 1171.5555
$ws.Range("I107").Value = 1130.3334
$ws.Range("J107").Value = 1192.1666
$ws.Range("K107").Value = 1130.3334
$ws.Range("L107").Value = 1192.1666
$ws.Range("M107").Value = 789.6666
$ws.Range("N107").Value = -5032.1666
# Row 122
$ws.Range("H122").Value = 1493
$ws.Range("I122").Value = 1217.5714
$ws.Range("K122").Value = 3652.7142
$ws.Range("M122").Value = -1202.7142
# Row 134
$ws.Range("H134").Value = 1995.1724
$ws.Range("I134").Value = 1145.3334
$ws.Range("J134").Value = 4226
$ws.Range("K134").Value = 3436.0002
$ws.Range("L134").Value = 12678
$ws.Range("M134").Value = -901.0001999999999
$ws.Range("N134").Value = -17748
# Row 136
$ws.Range("H136").Value = 1011.1316
$ws.Range("I136").Value = 796.5333000000001
$ws.Range("J136").Value = 1815.875
$ws.Range("K136").Value = 2389.5999
$ws.Range("L136").Value = 5447.625
$ws.Range("M136").Value = 160.4000999999998
$ws.Range("N136").Value = -10547.625

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 5.1666665
$ws.Range("I12").Value = 11
$ws.Range("J12").Value = 2.9230769
$ws.Range("K12").Value = 33
$ws.Range("L12").Value = 8.7692307
$ws.Range("M12").Value = 140
$ws.Range("N12").Value = -354.7692307
# Row 80
$ws.Range("H80").Value = 2000
$ws.Range("J80").Value = 2000
$ws.Range("L80").Value = 6000
$ws.Range("N80").Value = -7872
# Row 83
$ws.Range("H83").Value = 2000
$ws.Range("J83").Value = 2000
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27360
# Row 113
$ws.Range("H113").Value = 663787
$ws.Range("I113").Value = 1078218.4
$ws.Range("J113").Value = 696.9
$ws.Range("K113").Value = 3234655.2
$ws.Range("L113").Value = 2090.7
$ws.Range("M113").Value = -3232485.2
$ws.Range("N113").Value = -6430.7
# Row 132
$ws.Range("H132").Value = 920035.0600000001
$ws.Range("J132").Value = 920035.0600000001
$ws.Range("L132").Value = 8280315.540000001
$ws.Range("N132").Value = -8285375.540000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 6984.5557
$ws.Range("I113").Value = 932.75
$ws.Range("J113").Value = 11826
$ws.Range("K113").Value = 932.75
$ws.Range("L113").Value = 11826
$ws.Range("M113").Value = 1237.25
$ws.Range("N113").Value = -16166

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 572.6
$ws.Range("I22").Value = 576.4286
$ws.Range("J22").Value = 563.6667
$ws.Range("K22").Value = 576.4286
$ws.Range("L22").Value = 563.6667
$ws.Range("M22").Value = -281.4286
$ws.Range("N22").Value = -1153.6667
# Row 27
$ws.Range("H27").Value = 572.6
$ws.Range("I27").Value = 576.4286
$ws.Range("J27").Value = 563.6667
$ws.Range("K27").Value = 576.4286
$ws.Range("L27").Value = 563.6667
$ws.Range("M27").Value = -469.4286
$ws.Range("N27").Value = -777.6667
# Row 46
$ws.Range("H46").Value = 666.5
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 733
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 733
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -1109
# Row 50
$ws.Range("H50").Value = 10000
$ws.Range("J50").Value = 10000
$ws.Range("L50").Value = 10000
$ws.Range("N50").Value = -11274
# Row 55
$ws.Range("H55").Value = 410.8095
$ws.Range("I55").Value = 468.91666
$ws.Range("J55").Value = 333.33334
$ws.Range("K55").Value = 468.91666
$ws.Range("L55").Value = 333.33334
$ws.Range("M55").Value = -295.91666
$ws.Range("N55").Value = -679.33334
# Row 56
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = ""
# Row 122
$ws.Range("H122").Value = 2980
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2980
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8940
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -13840

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 14175
$ws.Range("J92").Value = 14175
$ws.Range("L92").Value = 14175
$ws.Range("N92").Value = -19167
# Row 100
$ws.Range("H100").Value = 636.7619
$ws.Range("I100").Value = 342.63635
$ws.Range("J100").Value = 960.3
$ws.Range("K100").Value = 685.2727
$ws.Range("L100").Value = 1920.6
$ws.Range("M100").Value = -144.2727
$ws.Range("N100").Value = -3002.6
# Row 122
$ws.Range("H122").Value = 3286.111
$ws.Range("I122").Value = 3448.3333
$ws.Range("J122").Value = 2961.6667
$ws.Range("K122").Value = 10344.9999
$ws.Range("L122").Value = 8885.000100000001
$ws.Range("M122").Value = -7894.999899999999
$ws.Range("N122").Value = -13785.0001
# Row 132
$ws.Range("H132").Value = 1940.5161
$ws.Range("I132").Value = 1352.1
$ws.Range("J132").Value = 3010.3635
$ws.Range("K132").Value = 4056.3
$ws.Range("L132").Value = 9031.0905
$ws.Range("M132").Value = -1526.3
$ws.Range("N132").Value = -14091.0905

